# Apply updated figures to the "Fonds de solidarite" (VOLET1) worksheet.
# For a set of rows, the nombre_aides (col C) and montant_total (col E)
# values are refreshed with newer totals as of the 2022-06-03 data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new C value, new E value)
$updates = @{
    2   = @(100826, 327347554)
    63  = @(14361,  36190509)
    74  = @(950,    4266554)
    83  = @(3412,   115697749)
    92  = @(409127, 1595300380)
    93  = @(209592, 1309153481)
    94  = @(94195,  917486940)
    95  = @(50768,  932648032)
    96  = @(17287,  793812185)
    99  = @(165,    25935853)
    100 = @(413,    66473777)
    104 = @(135242, 272218580)
    175 = @(80783,  486175817)
    177 = @(14718,  251572342)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
